{"js": "// Turn the old \"Heading1\" title + bold \"By <Author>\" byline into the\n// pandoc-style title-block pair: a \"Title\"-styled heading paragraph and a\n// separate \"Authors\"-styled paragraph holding just the author name(s).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Title paragraph: switch from Heading1 to the Title style. -----------\nconst titlePara = paragraphs.items[0];\ntitlePara.style = \"Title\";\n\n// --- Byline paragraph: \"By Dorothy Day\" (bold) -> \"Dorothy Day\" (Authors) -\nconst authorPara = paragraphs.items[1];\nauthorPara.clear();\nauthorPara.insertText(\"Dorothy Day\", \"Start\");\nauthorPara.style = \"Authors\";\n\nawait context.sync();\n", "ps1": "# Turn the old \"Heading1\" title + bold \"By <Author>\" byline into the\n# pandoc-style title-block pair: a \"Title\"-styled heading paragraph and a\n# separate \"Authors\"-styled paragraph holding just the author name(s).\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: switch from Heading1 to the Title style. -----------\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Style = \"Title\"\n\n# --- Byline paragraph: \"By Dorothy Day\" (bold) -> \"Dorothy Day\" (Authors) -\n$authorPara = $d.Paragraphs.Item(2)\n# Drop the bold run formatting before we touch the text/style so the new\n# run comes out with no leftover direct character formatting.\n$authorPara.Range.Font.Bold = 0\n$authorPara.Range.Text = \"Dorothy Day\"\n$authorPara.Style = \"Authors\"\n"}
